$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Usuario changes from "apellegrini" to "ocerutti"
$ws.Range("C3").Value = "ocerutti"

# Row 4 gains its own Usuario / Contrasenia values
$ws.Range("C4").Value = "apellegrini"
$ws.Range("D4").Value = "silverarrow"

# NroSiniestro value in F3 increments by 1 (leading apostrophe keeps it text, matching the existing quotePrefix text format)
$ws.Range("F3").Value = "'1120170200908"

# Update the selected cell to D8
$ws.Range("D8").Select()
